# "Generate Report for Handback" - refresh the CI handback-status report
# with the latest handoff/handback timestamps for the file whose row sits
# at row 2 of each per-locale worksheet (the other tracked file, row 3,
# keeps the timestamps from its last recorded run).
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 16:59:19"   # Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-03-24 16:59:45"   # Correspond Handback DateTime

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 16:59:24"   # Correspond Handoff Datetime
$wsDeDe.Range("H2").Value = "2016-03-24 16:59:51"   # Correspond Handback DateTime
